# Update countries & provincias Spain
# Refresh the "Pais" COVID stats sheet: bump the "last updated" timestamp,
# push new totals for several countries, and swap a handful of country
# labels whose rank changed so the sheet stays in its sorted order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row([int]$r, $country, $b, $c, $d, $e, $f, $g, $h) {
    if ($country -ne $null) {
        $ws.Cells.Item($r, 1).Value2 = $country
    }
    $ws.Cells.Item($r, 2).Value2 = $b
    $ws.Cells.Item($r, 3).Value2 = $c
    $ws.Cells.Item($r, 4).Value2 = $d
    $ws.Cells.Item($r, 5).Value2 = $e
    $ws.Cells.Item($r, 6).Value2 = $f
    $ws.Cells.Item($r, 7).Value2 = $g
    $ws.Cells.Item($r, 8).Value2 = $h
}

# --- Header: "Datos actualizados ..." timestamp (A1) -----------------------
$ws.Range("A1").Value2 = "Datos actualizados a 4 de Octubre de 2020 a las 11:13"

# --- Rows whose data refreshed but whose country (column A) is unchanged ---
Set-Row 5  $null 6553027  5614 5509966 941220 0  29 101841
Set-Row 7  $null 1215001 10499  979143 214500 0 107  21358
Set-Row 23 $null  322497  3190  273079  43642 0 100   5776
Set-Row 27 $null  265086   643  193027  70367 0  10   1692
Set-Row 57 $null   72310     0   66813   5238 0   1    259
Set-Row 65 $null   48146   714   38629   8704 0   4    813
Set-Row 89 $null   17659   258   15849   1512 0   5    298
Set-Row 98 $null   13139   818    4828   8256 0   1     55

# --- Oman / Polonia swap places (rows 42-43) --------------------------------
Set-Row 42 "Polonia" 100074 1934 73003 24441 0 26 2630
Set-Row 43 "Oman"      98585    0 88528  9122 0  0  935

# --- Maldivas / Finlandia swap places (rows 104-105) ------------------------
Set-Row 104 "Finlandia" 10538 147 8100 2093 0 0 345
Set-Row 105 "Maldivas"  10465   0 9310 1121 0 0  34

# --- Nicaragua / Hong Kong / Congo / Lituania reshuffle (rows 125-128) ------
Set-Row 125 "Lituania"  5185 104 2495 2596 0 0  94
Set-Row 126 "Nicaragua" 5170   0 2913 2106 0 0 151
Set-Row 127 "Hong Kong" 5109   0 4849  155 0 0 105
Set-Row 128 "Congo"     5089   0 3887 1113 0 0  89

# --- Somalia / Gambia / Tailandia / Estonia reshuffle (rows 139-142) --------
Set-Row 139 "Estonia"   3603  26 2749  787 0 0  67
Set-Row 140 "Somalia"   3593   0 3001  493 0 0  99
Set-Row 141 "Gambia"    3590   0 2226 1249 0 0 115
Set-Row 142 "Tailandia" 3585   2 3388  138 0 0  59
